$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) and Volume (E) columns so numeric-looking
# strings (e.g. "1.00", "0.998") are preserved exactly instead of being
# auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "66.812.18"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "2.476.83"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "583.52"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").Value = "167.20"
$ws.Range("E6").Value = "  -5.64%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E8").Value = "  -3.03%  "
$ws.Range("D9").Value = "2.480.29"
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("D10").Value = "0.135"
$ws.Range("E10").Value = "  -4.27%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "0.336"
$ws.Range("E12").Value = "  -3.20%  "
$ws.Range("D13").Value = "4.93"
$ws.Range("E13").Value = "  -4.25%  "
$ws.Range("D14").Value = "25.81"
$ws.Range("E14").Value = "  -4.31%  "
$ws.Range("D15").Value = "2.920.88"
$ws.Range("E15").Value = "  -2.34%  "
$ws.Range("D16").Value = "0.0000173"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D17").Value = "66.780.93"
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("D18").Value = "2.476.39"
$ws.Range("E18").Value = "  -2.22%  "
$ws.Range("D19").Value = "11.60"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").Value = "7.86"
$ws.Range("E20").Value = "  -1.88%  "
$ws.Range("D21").Value = "362.03"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("E22").Value = "  -4.23%  "
$ws.Range("D23").Value = "4.45"
$ws.Range("E23").Value = "  -5.45%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "70.82"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("E26").Value = "  -6.13%  "
$ws.Range("D27").Value = "9.44"
$ws.Range("E27").Value = "  -7.95%  "
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("D30").Value = "0.0₃0930"
$ws.Range("E30").Value = "  -6.29%  "
$ws.Range("D31").Value = "8.05"
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("D32").Value = "512.68"
$ws.Range("E32").Value = "  -5.68%  "
$ws.Range("E33").Value = "  -2.65%  "
$ws.Range("D34").Value = "1.27"
$ws.Range("E34").Value = "  -5.64%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Value = "0.126"
$ws.Range("E36").Value = "  -2.42%  "
$ws.Range("D37").Value = "158.38"
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("E38").Value = "  -3.64%  "
$ws.Range("D39").Value = "18.91"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").Value = "18.55"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("E41").Value = "  -3.61%  "
$ws.Range("D42").Value = "4.94"
$ws.Range("E42").Value = "  -5.34%  "
$ws.Range("D43").Value = "0.333"
$ws.Range("E43").Value = "  -6.77%  "
$ws.Range("D44").Value = "2.44"
$ws.Range("E44").Value = "  -3.24%  "
$ws.Range("D45").Value = "39.25"
$ws.Range("E45").Value = "  -2.00%  "
$ws.Range("D46").Value = "142.29"
$ws.Range("E46").Value = "  -3.36%  "
$ws.Range("D47").Value = "0.538"
$ws.Range("E47").Value = "  -4.37%  "
$ws.Range("D48").Value = "3.58"
$ws.Range("E48").Value = "  -4.02%  "
$ws.Range("E49").Value = "  -3.81%  "
$ws.Range("D50").Value = "1.64"
$ws.Range("E50").Value = "  -3.65%  "
$ws.Range("D51").Value = "0.0737"
$ws.Range("E51").Value = "  -2.69%  "

# Restore default (Normal) style so no stray number-format style is left
# behind now that the text values are set.
$ws.Range("D2:E51").Style = "Normal"

